$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    if ($val -match "^-?[0-9]*\.?[0-9]+$") {
        $c.NumberFormat = "@"
        $c.Value = $val
        $c.Style = "Normal"
    } else {
        $c.Value = $val
    }
}

Set-TextValue 'D2' '27.717.89'
Set-TextValue 'E2' '  +1.64%  '
Set-TextValue 'D3' '1.881.44'
Set-TextValue 'E3' '  +1.65%  '
Set-TextValue 'E4' '  +0.49%  '
Set-TextValue 'D5' '333.39'
Set-TextValue 'E5' '  +2.70%  '
Set-TextValue 'D6' '1.007'
Set-TextValue 'E6' '  +0.44%  '
Set-TextValue 'D7' '0.4708'
Set-TextValue 'E7' '  +3.63%  '
Set-TextValue 'D8' '0.3945'
Set-TextValue 'E8' '  +1.78%  '
Set-TextValue 'D9' '47.97'
Set-TextValue 'E9' '  -0.70%  '
Set-TextValue 'D10' '0.08071'
Set-TextValue 'E10' '  +2.08%  '
Set-TextValue 'D11' '1.028'
Set-TextValue 'E11' '  +1.56%  '
Set-TextValue 'D12' '22.18'
Set-TextValue 'E12' '  +3.97%  '
Set-TextValue 'D13' '1.881.59'
Set-TextValue 'E13' '  +0.66%  '
Set-TextValue 'D14' '5.978'
Set-TextValue 'E14' '  +1.38%  '
Set-TextValue 'D15' '7.132'
Set-TextValue 'E15' '  +0.02%  '
Set-TextValue 'D16' '1.008'
Set-TextValue 'E16' '  +0.54%  '
Set-TextValue 'D17' '87.21'
Set-TextValue 'E17' '  +1.76%  '
Set-TextValue 'D18' '0.06702'
Set-TextValue 'E18' '  +1.37%  '
Set-TextValue 'D19' '0.00001045'
Set-TextValue 'E19' '  +1.89%  '
Set-TextValue 'D20' '17.34'
Set-TextValue 'E20' '  +0.96%  '
Set-TextValue 'D21' '1.006'
Set-TextValue 'E21' '  +0.50%  '
Set-TextValue 'B22' 'Uniswap'
Set-TextValue 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D22' '5.538'
Set-TextValue 'E22' '  +1.05%  '
Set-TextValue 'B23' 'WrappedBTC'
Set-TextValue 'C23' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D23' '27.712.89'
Set-TextValue 'E23' '  +1.56%  '
Set-TextValue 'D24' '11.03'
Set-TextValue 'E24' '  +1.93%  '
Set-TextValue 'D25' '2.312'
Set-TextValue 'E25' '  +0.54%  '
Set-TextValue 'D26' '2.103.89'
Set-TextValue 'E26' '  +1.08%  '
Set-TextValue 'D27' '159.88'
Set-TextValue 'E27' '  +3.70%  '
Set-TextValue 'D28' '20.23'
Set-TextValue 'E28' '  +1.78%  '
Set-TextValue 'D29' '2.105'
Set-TextValue 'E29' '  +2.63%  '
Set-TextValue 'D30' '5.570'
Set-TextValue 'E30' '  +2.59%  '
Set-TextValue 'D31' '121.98'
Set-TextValue 'E31' '  +0.73%  '
Set-TextValue 'D32' '0.9828'
Set-TextValue 'E32' '  +4.68%  '
Set-TextValue 'D33' '0.09492'
Set-TextValue 'E33' '  +1.84%  '
Set-TextValue 'D34' '1.449'
Set-TextValue 'E34' '  +0.34%  '
Set-TextValue 'D35' '3.610'
Set-TextValue 'E35' '  +0.62%  '
Set-TextValue 'D36' '5.360'
Set-TextValue 'E36' '  +2.18%  '
Set-TextValue 'D37' '0.06134'
Set-TextValue 'E37' '  +1.76%  '
Set-TextValue 'D38' '0.02262'
Set-TextValue 'E38' '  +1.86%  '
Set-TextValue 'D39' '1.229'
Set-TextValue 'E39' '  +2.04%  '
Set-TextValue 'D40' '8.133'
Set-TextValue 'E40' '  +1.06%  '
Set-TextValue 'D41' '0.5998'
Set-TextValue 'E41' '  +1.51%  '
Set-TextValue 'D42' '0.1898'
Set-TextValue 'E42' '  +0.93%  '
Set-TextValue 'D43' '10.32'
Set-TextValue 'E43' '  +1.87%  '
Set-TextValue 'D44' '1.262'
Set-TextValue 'E44' '  -0.55%  '
Set-TextValue 'D45' '0.5727'
Set-TextValue 'E45' '  +2.30%  '
Set-TextValue 'E46' '  +0.67%  '
Set-TextValue 'B47' 'NEARProtocol'
Set-TextValue 'C47' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D47' '1.949'
Set-TextValue 'E47' '  +2.16%  '
Set-TextValue 'B48' 'PancakeSwap'
Set-TextValue 'C48' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D48' '3.395'
Set-TextValue 'E48' '  +0.63%  '
Set-TextValue 'D49' '0.06926'
Set-TextValue 'E49' '  +2.92%  '
Set-TextValue 'D50' '114.68'
Set-TextValue 'E50' '  +6.51%  '
Set-TextValue 'D51' '1.073'
Set-TextValue 'E51' '  +2.14%  '
